# Applies the "Updated cryptos list" data refresh (GitHub Actions commit)
# to the cryptocurrency table on Sheet1 (columns B-E, rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.803.74"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "2.217.60"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.13"
$ws.Range("E5").Value = "  +5.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.630"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.41"
$ws.Range("E7").Value = "  +2.48%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.601"
$ws.Range("E9").Value = "  +8.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.52"
$ws.Range("E10").Value = "  +10.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0963"
$ws.Range("E11").Value = "  -3.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.31"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.19"
$ws.Range("E13").Value = "  +6.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.105"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("D15").Value = "2.551.29"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.94"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "2.230.30"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").Value = "41.708.72"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("D20").Value = "0.0₃0961"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.20"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.75"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.43"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +3.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.00"
$ws.Range("E25").Value = "  +8.89%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.53"
$ws.Range("E27").Value = "  +4.77%  "
$ws.Range("E28").Value = "  +8.46%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "170.60"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.70"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.55"
$ws.Range("E33").Value = "  +4.39%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.124"
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0737"
$ws.Range("E35").Value = "  +2.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.72"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.38"
$ws.Range("E37").Value = "  +16.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.01"
$ws.Range("E38").Value = "  +7.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0299"
$ws.Range("E39").Value = "  +8.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.28"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.94"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.36"
$ws.Range("E42").Value = "  +20.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.77"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.205"
$ws.Range("E44").Value = "  +6.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.83"
$ws.Range("E45").Value = "  -2.26%  "
$ws.Range("B46").Value = "SynthetixNetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.74"
$ws.Range("E46").Value = "  +2.80%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.68"
$ws.Range("E47").Value = "  -6.90%  "
$ws.Range("E48").Value = "  -1.82%  "
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("E50").Value = "  +5.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.40"
$ws.Range("E51").Value = "  +2.61%  "
